$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column D header, matching the style of the existing header cells (A1:C1)
$ws.Range("D1").Value = "Templates"
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update data rows: new date serials, updated "Sent at" timestamps, and new template file names
$ws.Range("B2").Value = 46020
$ws.Range("C2").Value = "Sent at 2025-12-29 11:59:21"
$ws.Range("D2").Value = "Customer Mail New.msg"

$ws.Range("B3").Value = 46020
$ws.Range("C3").Value = "Sent at 2025-12-29 11:59:23"
$ws.Range("D3").Value = "Power BI.msg"
